# Adds a new "docIndex" column (B) holding the zero-based row index of each
# record, shifting the existing "contents"/"source" columns one place right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of existing data rows (including the header row) before we touch
# anything, so the fill loop below still works if the sheet size changes.
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new blank column before column B ("contents" becomes C,
# "source" becomes D).
[void]$ws.Columns("B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "docIndex"

# Fill the new column with the 0-based document/cluster index for each
# data row (row 2 => index 0, row 3 => index 1, ...).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $r - 2
}

# Restore the active selection like the authored workbook.
[void]$ws.Range("D22").Select()
